$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value  = -22.41030000000002
$ws.Range("A10").Value = -22.03499999999999
$ws.Range("A12").Value = -21.4097
$ws.Range("A18").Value = -22.20140000000002
$ws.Range("A25").Value = -21.72529999999999
$ws.Range("A37").Value = -19.81329999999999
$ws.Range("A55").Value = -22.1871
$ws.Range("A68").Value = -21.46499999999999
$ws.Range("A77").Value = -20.12809999999999
$ws.Range("A78").Value = -19.69169999999998
$ws.Range("A79").Value = -20.27009999999999
$ws.Range("A80").Value = -19.3276
$ws.Range("A81").Value = -22.1226
$ws.Range("A82").Value = -21.63210000000001
$ws.Range("A84").Value = -21.9737
$ws.Range("A101").Value = -20.62519999999998
$ws.Range("A102").Value = -20.32619999999998
